$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the extra-training entry: "Dr. Brad Boehmke" -> "Dr. Bradley Boehmke"
$ws.Range("E11").Value = "Two-day workshop lead by Dr. Bradley Boehmke"

# Remove the "Reading Research Quarterly" entry from the peer-reviewer list
# (row 26 in the original sheet). Deleting the whole row shifts rows 27-37
# up to become 26-36, matching the new A1:H36 dimension, and also drops the
# now-unused shared string automatically.
$ws.Rows(26).Delete()

# Update the active selection to E11, matching the saved workbook state
$ws.Range("E11").Select()
